$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove rows 4-7 (old extra glycan entries no longer present)
$ws.Range("A4:J7").Delete() | Out-Null

# Insert two new columns before old column D (flexibility) to host sum_SASA and max_SASA
$ws.Range("D1:E1").EntireColumn.Insert() | Out-Null

# --- Header row (row 1) ---
$ws.Range("A1").Value = "glycan"
$ws.Range("B1").Value = "binding_score"
$ws.Range("C1").Value = "SASA"
$ws.Range("D1").Value = "sum_SASA"
$ws.Range("E1").Value = "max_SASA"
$ws.Range("F1").Value = "flexibility"
$ws.Range("G1").Value = "Q"
$ws.Range("H1").Value = "theta"
$ws.Range("I1").Value = "conformation"
$ws.Range("J1").Value = "monosaccharides"
$ws.Range("K1").Value = "motifs"
$ws.Range("L1").Value = "class"

# --- Row 2 data ---
$ws.Range("A2").Value = "Gal(b1-4)GlcNAc(b1-2)Man(a1-3)[Gal(b1-4)GlcNAc(b1-2)Man(a1-6)]Man(b1-4)GlcNAc(b1-4)GlcNAc"
$ws.Range("B2").Value = 1.205755702476461
$ws.Range("C2").Value = 1.901990340784573
$ws.Range("D2").Value = 13.31393238549201
$ws.Range("E2").Value = 2.596661214730085
$ws.Range("F2").Value = 24.47100952380953
$ws.Range("G2").Value = 0.4384285714285714
$ws.Range("H2").Value = 10.41
$ws.Range("I2").Value = "['4C1', '4C1', '4C1', '4C1', '4C1', '4C1', '4C1']"
$ws.Range("J2").Value = "['Gal(b1-4)', 'GlcNAc(b1-2)', 'Man(a1-3)', 'Man(b1-4)', 'Man(a1-6)', 'GlcNAc(b1-2)', 'Gal(b1-4)']"
$ws.Range("K2").Value = "['Gal(b1-3/4)GlcNAc(b1-2)Man(a1-3)[Gal(b1-3/4)GlcNAc(b1-2)Man(a1-6)]Man']"
$ws.Range("L2").Value = "N"

# --- Row 3 data ---
$ws.Range("A3").Value = "Neu5Ac(a2-3)Gal(b1-4)GlcNAc(b1-2)Man(a1-3)[Neu5Ac(a2-3)Gal(b1-4)GlcNAc(b1-2)Man(a1-6)]Man(b1-4)GlcNAc(b1-4)GlcNAc"
$ws.Range("B3").Value = -0.1785431686174741
$ws.Range("C3").Value = 1.682105456660559
$ws.Range("D3").Value = 11.77473819662391
$ws.Range("E3").Value = 2.425021331076254
$ws.Range("F3").Value = 27.30604857142857
$ws.Range("G3").Value = 0.4505714285714286
$ws.Range("H3").Value = 6.141428571428571
$ws.Range("I3").Value = "['4C1', '4C1', '4C1', '4C1', '4C1', '4C1', '4C1']"
$ws.Range("J3").Value = "['Gal(b1-4)', 'GlcNAc(b1-2)', 'Man(a1-3)', 'Man(b1-4)', 'Man(a1-6)', 'GlcNAc(b1-2)', 'Gal(b1-4)']"
$ws.Range("K3").Value = "['Gal(b1-3/4)GlcNAc(b1-2)Man(a1-3)[Gal(b1-3/4)GlcNAc(b1-2)Man(a1-6)]Man']"
$ws.Range("L3").Value = "N"

Write-Host "Edit applied"
